$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 325370
$ws.Range("D2").Value = 414484961
$ws.Range("C10").Value = 118238
$ws.Range("D10").Value = 173244238
$ws.Range("C11").Value = 144
$ws.Range("D11").Value = 213013
$ws.Range("C12").Value = 60313
$ws.Range("D12").Value = 87042626
$ws.Range("C16").Value = 4032
$ws.Range("D16").Value = 5722697
$ws.Range("C20").Value = 6919
$ws.Range("D20").Value = 9658932
$ws.Range("C22").Value = 78403
$ws.Range("D22").Value = 97695968
$ws.Range("C28").Value = 32744
$ws.Range("D28").Value = 47926911
$ws.Range("C30").Value = 11612
$ws.Range("D30").Value = 16703214
$ws.Range("C35").Value = 1901
$ws.Range("D35").Value = 2684626
$ws.Range("C36").Value = 98315
$ws.Range("D36").Value = 123653914
$ws.Range("C44").Value = 44728
$ws.Range("D44").Value = 65548733
$ws.Range("C46").Value = 9269
$ws.Range("D46").Value = 13294085
$ws.Range("C51").Value = 2417
$ws.Range("D51").Value = 3377321
$ws.Range("C52").Value = 69882
$ws.Range("D52").Value = 87639529
$ws.Range("C57").Value = 385
$ws.Range("D57").Value = 565730
$ws.Range("C59").Value = 28463
$ws.Range("D59").Value = 41740846
$ws.Range("C62").Value = 11279
$ws.Range("D62").Value = 16309419
$ws.Range("C64").Value = 1368
$ws.Range("D64").Value = 1912237
$ws.Range("C68").Value = 1534
$ws.Range("D68").Value = 2147591
$ws.Range("C70").Value = 20739
$ws.Range("D70").Value = 27157590
$ws.Range("C74").Value = 7672
$ws.Range("D74").Value = 11234856
$ws.Range("C76").Value = 5187
$ws.Range("D76").Value = 7531486
$ws.Range("C77").Value = 495
$ws.Range("D77").Value = 701239
$ws.Range("C78").Value = 288
$ws.Range("D78").Value = 405083
$ws.Range("C79").Value = 142679
$ws.Range("D79").Value = 177807291
$ws.Range("C85").Value = 64194
$ws.Range("D85").Value = 94080698
$ws.Range("C88").Value = 30106
$ws.Range("D88").Value = 43549256
$ws.Range("C90").Value = 2757
$ws.Range("D90").Value = 3969608
$ws.Range("C91").Value = 2932
$ws.Range("D91").Value = 4146126
$ws.Range("C92").Value = 33925
$ws.Range("D92").Value = 45982281
$ws.Range("C95").Value = 31
$ws.Range("D95").Value = 44814
$ws.Range("C96").Value = 8226
$ws.Range("D96").Value = 12091517
$ws.Range("C98").Value = 7571
$ws.Range("D98").Value = 10984190
$ws.Range("C101").Value = 504
$ws.Range("D101").Value = 727391
$ws.Range("C102").Value = 10692
$ws.Range("D102").Value = 16480189
$ws.Range("C104").Value = 2627
$ws.Range("D104").Value = 4329376
$ws.Range("C106").Value = 3578
$ws.Range("D106").Value = 5911207
$ws.Range("C110").Value = 143319
$ws.Range("D110").Value = 177218080
$ws.Range("C116").Value = 53237
$ws.Range("D116").Value = 78024109
$ws.Range("C118").Value = 27520
$ws.Range("D118").Value = 39870822
$ws.Range("C122").Value = 2331
$ws.Range("D122").Value = 3275199
$ws.Range("C124").Value = 524648
$ws.Range("D124").Value = 692891127
$ws.Range("C129").Value = 1394
$ws.Range("D129").Value = 2066182
$ws.Range("C131").Value = 211501
$ws.Range("D131").Value = 310905529
$ws.Range("C132").Value = 413
$ws.Range("D132").Value = 616250
$ws.Range("C134").Value = 187206
$ws.Range("D134").Value = 272242172
$ws.Range("C137").Value = 2875
$ws.Range("D137").Value = 4036252
$ws.Range("C139").Value = 6611
$ws.Range("D139").Value = 9342137
$ws.Range("C142").Value = 45474
$ws.Range("D142").Value = 60698389
$ws.Range("C148").Value = 14260
$ws.Range("D148").Value = 20907241
$ws.Range("C149").Value = 3819
$ws.Range("D149").Value = 5507100
$ws.Range("C152").Value = 405
$ws.Range("D152").Value = 582716
$ws.Range("C154").Value = 399
$ws.Range("D154").Value = 564263
$ws.Range("C155").Value = 17938
$ws.Range("D155").Value = 23708075
$ws.Range("C159").Value = 7298
$ws.Range("D159").Value = 10618344
$ws.Range("C161").Value = 5099
$ws.Range("D161").Value = 7340996
$ws.Range("C164").Value = 274
$ws.Range("D164").Value = 392364
$ws.Range("C166").Value = 18987
$ws.Range("D166").Value = 31173397
$ws.Range("C167").Value = 2067
$ws.Range("D167").Value = 3417822
$ws.Range("C171").Value = 109
$ws.Range("D171").Value = 188949
$ws.Range("C172").Value = 88870
$ws.Range("D172").Value = 111051119
$ws.Range("C179").Value = 34189
$ws.Range("D179").Value = 50133334
$ws.Range("C181").Value = 13199
$ws.Range("D181").Value = 19070357
$ws.Range("C183").Value = 1255
$ws.Range("D183").Value = 1757339
$ws.Range("C185").Value = 1712
$ws.Range("D185").Value = 2402871
$ws.Range("C187").Value = 241383
$ws.Range("D187").Value = 299920686
$ws.Range("C195").Value = 87371
$ws.Range("D195").Value = 128070825
$ws.Range("C198").Value = 33455
$ws.Range("D198").Value = 48160920
$ws.Range("C201").Value = 5163
$ws.Range("D201").Value = 7349896
$ws.Range("C204").Value = 5056
$ws.Range("D204").Value = 7000250
$ws.Range("C207").Value = 267489
$ws.Range("D207").Value = 330968798
$ws.Range("C216").Value = 96127
$ws.Range("D216").Value = 140622357
$ws.Range("C219").Value = 52192
$ws.Range("D219").Value = 75432161
$ws.Range("C222").Value = 4717
$ws.Range("D222").Value = 6620548
$ws.Range("C225").Value = 6008
$ws.Range("D225").Value = 8316779
$ws.Range("C228").Value = 107724
$ws.Range("D228").Value = 134651820
$ws.Range("C235").Value = 49978
$ws.Range("D235").Value = 73215137
$ws.Range("C237").Value = 12632
$ws.Range("D237").Value = 18163490
$ws.Range("C238").Value = 13
$ws.Range("D238").Value = 13183
$ws.Range("C241").Value = 2600
$ws.Range("D241").Value = 3640192
$ws.Range("C242").Value = 261549
$ws.Range("D242").Value = 330197014
$ws.Range("C243").Value = 174
$ws.Range("D243").Value = 215181
$ws.Range("C250").Value = 96860
$ws.Range("D250").Value = 141927087
$ws.Range("C253").Value = 66155
$ws.Range("D253").Value = 95873681
$ws.Range("C255").Value = 2437
$ws.Range("D255").Value = 3436773
$ws.Range("C257").Value = 7
$ws.Range("D257").Value = 10500
$ws.Range("C258").Value = 4787
$ws.Range("D258").Value = 6726683
